$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The vehicle profile data (distance / parking-at-home) was split out into
# separate uploaded files. The sheet no longer needs the placeholder
# "default" vehicle row (ID_Vehicle = 0) nor the stray duplicate trailing
# row (ID_Vehicle = 1994101), and every remaining vehicle now carries an
# explicit reference to its own parking-at-home and distance profile.

# Remove the placeholder default row (old row 2, ID_Vehicle = 0).
$ws.Rows(2).Delete()

# Remove the duplicate trailing row (now row 11, ID_Vehicle = 1994101).
$ws.Rows(11).Delete()

# Add the two new profile-id columns.
$ws.Range("N1").Value = "id_parking_at_home_profile"
$ws.Range("O1").Value = "id_distance_profile"

# Each vehicle's profile ids mirror its own ID_Vehicle value.
for ($r = 2; $r -le 10; $r++) {
    $idVal = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 14).Value = $idVal
    $ws.Cells.Item($r, 15).Value = $idVal
}

# Size the two new columns to fit their contents, matching the rest of the
# (best-fit) table.
$ws.Columns.Item(14).ColumnWidth = 21.92
$ws.Columns.Item(15).ColumnWidth = 14.76

# Restore the selection saved with the workbook.
$ws.Range("F20").Select()
